$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.817.49"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.66%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.629.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.74%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.17%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.21%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5072"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.24%  "

# Row 7
$ws.Range("E7").Value = "  +0.25%  "

# Row 8
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2581"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.04%  "

# Row 9
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06460"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.61%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.44"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.12%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07808"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.68%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.259"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.50%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.627.30"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.93%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.854.00"
$ws.Range("D14").Style = "Normal"

# Row 15
$ws.Range("E15").Value = "  +1.48%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.19"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.82%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0₅7541"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.63%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.823.90"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.74%  "

# Row 19
$ws.Range("E19").Value = "  +0.16%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "194.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.94%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.296"
$ws.Range("D21").Style = "Normal"

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.817"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.33%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.006"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.73%  "

# Row 24
$ws.Range("E24").Value = "  +0.07%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.815"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.17%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "140.45"
$ws.Range("D26").Style = "Normal"

# Row 27
$ws.Range("E27").Value = "  +2.05%  "

# Row 28
$ws.Range("E28").Value = "  -2.25%  "

# Row 29
$ws.Range("E29").Value = "  -1.41%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.237"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.08%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.04866"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.17%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.178"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.35%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.555"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.61%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.374"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.09%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.8942"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.28%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.573"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.04%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.132.23"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.70%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5463"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.42%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01558"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.08%  "

# Row 41
$ws.Range("E41").Value = "  -0.06%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.572"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.61%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7950"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.39%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "97.32"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.86%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.781.04"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.08%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₈114"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.49%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4436"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.22%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "55.09"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.78%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05057"
$ws.Range("D49").Style = "Normal"

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.626"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.26%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.004"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.08%  "

